$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: apply the built-in "Neutral" cell style to the "Tags" cell, and add two new
# to-do items alongside it. Write these first so the new shared-string entries land in
# the same order as the target workbook.
$ws.Range("A13").Style = "Neutral"
$ws.Range("C13").Value = "whitespace check"
$ws.Range("D13").Value = "add fileattachment support"

# New column F: copy the formatting already used by the header row (column E, style "Gut")
# onto F3 first, then set its text - this reuses the existing style record instead of
# creating a brand new duplicate "Good" style entry.
$ws.Range("E3").Copy($ws.Range("F3"))
$ws.Range("F3").Value = "change emote next to affinity to display affinity"

# Column F width (~27.71 characters); ColumnWidth only accepts 1/6-character increments
# here, so 26.9 lands on the closest achievable stored width.
$ws.Columns.Item(6).ColumnWidth = 26.9

# New row 22 item
$ws.Range("A22").Value = "Sentry service"

# Update the selected cell in the view
$ws.Range("A23").Select()
